$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The ground-names table contained two duplicate rows that need to be
# removed as part of cleaning up the training/checking data:
#   row 7 - "Melbourne / MCG / Melbourne Cricket Ground" (duplicate of row 5)
#   row 9 - "Optus Stadium / Perth Stadium" (duplicate of row 8)
# Delete from the bottom up so row numbers used for the later delete are
# not shifted by the earlier delete.
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(7).Delete()

# Re-sync the Ground_Id column (G) to the expected values for the
# remaining rows now that the two duplicate grounds are gone.
$ws.Range("G7").Value = 6
$ws.Range("G8").Value = 8
$ws.Range("G9").Value = 9
$ws.Range("G10").Value = 10
$ws.Range("G11").Value = 11
$ws.Range("G12").Value = 12
$ws.Range("G13").Value = 13
$ws.Range("G14").Value = 14
$ws.Range("G15").Value = 15
$ws.Range("G16").Value = 16
$ws.Range("G17").Value = 17
$ws.Range("G18").Value = 18
$ws.Range("G19").Value = 19
$ws.Range("G20").Value = 20
$ws.Range("G21").Value = 21
$ws.Range("G22").Value = 22
$ws.Range("G23").Value = 23
$ws.Range("G24").Value = 24
$ws.Range("G25").Value = 25
$ws.Range("G26").Value = 26

$ws.Range("A8").Select()
